# Applies the NMCARS-PART-5231 edit described by the commit:
#   "Fixed NMCARS internal xrefs and image sources for EMF files"
#
# Content-level changes (all inside the document body):
#   1. The "negotiations.   (g) Copies ..." paragraph switches from an
#      ad-hoc <w:widowControl w:val="0"/> paragraph property to the
#      "List1" paragraph style.
#   2. A run of previously un-styled "(n) ..." list paragraphs under
#      5231.205-90 (c) and (d) pick up the "List2" paragraph style.
#   3. Several "List3"-styled "(i)/(ii)/(iii) ..." paragraphs have their
#      leading "(i)"/"(ii)"/"(iii)" marker split out into its own run,
#      separate from the following text run.

$d = $word.ActiveDocument

# --- helper: split the run that begins paragraph $para after the first
# --- $len characters, without leaving any residual formatting markup.
# --- (Add a bookmark around the prefix, then immediately delete it --
# --- Word always breaks runs at bookmark boundaries, and removing the
# --- bookmark afterwards leaves the now-separate runs untouched.)
function Split-LeadingRun($para, [int]$len) {
    $full = $para.Range
    $prefix = $d.Range($full.Start, $full.Start + $len)
    $d.Bookmarks.Add("zzSplitMarker", $prefix) | Out-Null
    $d.Bookmarks("zzSplitMarker").Delete()
}

# 1. "negotiations.   (g) Copies ..." paragraph -> List1 style
#    (this also removes the old widowControl override, matching Word's
#    normal behaviour of swapping pPr contents when a style is applied)
$d.Paragraphs.Item(5).Style = "List1"

# 2. Un-styled "(n)" list paragraphs -> List2 style
$d.Paragraphs.Item(13).Style = "List2"   # (1) ASN(RDA) must make a determination...
$d.Paragraphs.Item(14).Style = "List2"   # (2) An agreement generally will be considered...
$d.Paragraphs.Item(15).Style = "List2"   # (3) The agreement shall apply to prospective...
$d.Paragraphs.Item(16).Style = "List2"   # (4) The agreement must project an overall benefit...
$d.Paragraphs.Item(18).Style = "List2"   # (1) The agreement shall require the contractor...
$d.Paragraphs.Item(22).Style = "List2"   # (2) The agreement shall require that the sum...
$d.Paragraphs.Item(23).Style = "List2"   # (3) The Navy may agree to modify the amount...

# 3. Split the leading "(i)"/"(ii)"/"(iii)" marker into its own run
Split-LeadingRun $d.Paragraphs.Item(19) 3    # "(i) The direct costs..."
Split-LeadingRun $d.Paragraphs.Item(20) 4    # "(ii) The incremental indirect costs..."
Split-LeadingRun $d.Paragraphs.Item(21) 5    # "(iii) The non-incremental indirect costs..."
Split-LeadingRun $d.Paragraphs.Item(24) 3    # "(i) Any smaller amount..."
Split-LeadingRun $d.Paragraphs.Item(25) 4    # "(ii) Any larger amount..."
Split-LeadingRun $d.Paragraphs.Item(26) 5    # "(iii) In determining whether..."
